$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 was an empty trailing row (A11 only carried the "key" style,
# no value). Duplicate its formatting down into row 12 first, so the
# blank trailing row keeps the same look after we repurpose row 11 for
# new data.
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill the previously-blank row 11 with the new key/value pair.
$ws.Range("A11").Value = "x_date_format"
$ws.Range("B11").Value = "month\nyyyy"
